# Update the "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" worksheets, matching the data
# refresh captured in the target diff.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 156
    3  = 7085
    4  = 4744
    5  = 70
    6  = 158
    8  = 41
    9  = 94
    10 = 68
    11 = 65
    14 = 128
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
